# Update the handoff status for the c9dbe6c4-a982-4f07-ab11-8edda0a0ec07.md
# file from "Ready for handoff" back to "In Translation" across all report
# sheets, as part of generating the status report for archive.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B4").Value = "In Translation"
$overview.Range("C4").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B4").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B4").Value = "In Translation"

$wb.Save()
